# Addressed small editorial changes.
# Change-Id: Ia5384b8439ddbe7865c4eaf0f997b7daf058fe8e

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute(
        $find,    # FindText
        $true,    # MatchCase
        $false,   # MatchWholeWord
        $false,   # MatchWildcards
        $false,   # MatchSoundsLike
        $false,   # MatchAllWordForms
        $true,    # Forward
        1,        # Wrap (wdFindContinue)
        $false,   # Format
        $replace, # ReplaceWith
        2)        # Replace (wdReplaceAll)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1. Fix duplicated "on on" typo.
Replace-Text `
    "Created Javascript framework (built on on top of Backbone.js)" `
    "Created Javascript framework (built on top of Backbone.js)"

# 2. Fix the displayed WiserTogether link text (https -> http). Using
# Hyperlink.TextToDisplay (rather than a blanket Find/Replace) updates just
# the visible text and keeps both the run's Link character style and the
# underlying hyperlink target/relationship untouched.
$fixed = $false
foreach ($h in $d.Hyperlinks) {
    if ($h.Address -eq "https://wisertogether.com") {
        $h.TextToDisplay = "http://wisertogether.com"
        $fixed = $true
        break
    }
}
if (-not $fixed) {
    throw "Could not find the wisertogether.com hyperlink"
}

# 3. Fix "20011" typo year.
Replace-Text "20011 - 2013" "2011 - 2013"

# 4. Pluralize "increase" -> "increases".
Replace-Text `
    "Oversaw and negotiated software staff increase as well as company merger" `
    "Oversaw and negotiated software staff increases as well as company merger"

# 5. Reword the first Textron Lycoming bullet.
Replace-Text `
    "Designed and provided ongoing development and maintenance of national award-winning dynamic Web site" `
    "Worked throughout high school and college, onsite and telecommuting, with the Information Technology Department"

# 6. Insert a new bulleted paragraph right before "Created a unique Intranet...",
#    carrying a reworded version of the bullet replaced in step 5, keeping the
#    same list numbering (numId 9 / ilvl 0).
#    NOTE: once InsertParagraphBefore() splits the paragraph, $target itself
#    becomes the new (now-empty) paragraph, while the original text shifts
#    into the following paragraph - so we just set $target's own text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Created a unique Intranet site providing executives*") {
        $target = $p
        break
    }
}
if ($null -eq $target) {
    throw "Could not find the 'Created a unique Intranet site' paragraph"
}
$target.Range.InsertParagraphBefore()
$target.Range.Text = "Designed and provided ongoing development and maintenance of company's dynamic website, which won a national design award"

# 7. Clarify the external marketing site bullet.
Replace-Text `
    "Created external marketing site, parts database" `
    "Created external marketing site for aircraft engine parts database"
